$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("drone")

# Row 2
$ws.Range("A2").Value = "Holy Stone HS170 Predator Mini RC Helicopter Drone 2.4Ghz 6-Axis Gyro 4 Channels Quadcopter Good Choice for Drone Training"
$ws.Range("B2").Value = "https://www.amazon.com/Holy-Stone-Predator-Helicopter-Quadcopter/dp/B0157IHJMQ/ref=sr_1_4?s=toys-and-games&ie=UTF8&qid=1528131930&sr=1-4&keywords=drone"
$ws.Range("C2").Value = "4.4 out of 5 stars"
$ws.Range("D2").Value = 4472
$ws.Range("D2").NumberFormat = "#,##0"
$ws.Range("E2").Value = 39

# Row 3
$ws.Range("A3").Value = "Holy Stone GPS FPV RC Drone HS100 with Camera Live Video and GPS Return Home Quadcopter with Adjustable Wide-Angle 720P HD WIFI Camera- Follow Me, Altitude Hold, Intelligent Battery Long Control Range"
$ws.Range("B3").Value = "https://www.amazon.com/Holy-Stone-Quadcopter-Adjustable-Intelligent/dp/B074YYVXQH/ref=sr_1_5?s=toys-and-games&ie=UTF8&qid=1528131930&sr=1-5&keywords=drone"
$ws.Range("C3").Value = "4.6 out of 5 stars"
$ws.Range("D3").Value = 930
$ws.Range("E3").Value = 279

# Row 4
$ws.Range("A4").Value = "Holy Stone HS160 Shadow FPV RC Drone with 720P HD Wi-Fi Camera Live Video Feed 2.4GHz 6-Axis Gyro Quadcopter for Kids & Beginners - Altitude Hold, One Key Start, Foldable Arms,Bonus Battery"
$ws.Range("B4").Value = "https://www.amazon.com/Holy-Stone-Shadow-Quadcopter-Beginners/dp/B074S2HK59/ref=sr_1_6?s=toys-and-games&ie=UTF8&qid=1528131930&sr=1-6&keywords=drone"
$ws.Range("C4").Value = "4.3 out of 5 stars"
$ws.Range("D4").Value = 1015
$ws.Range("E4").Value = 99

# Row 5
$ws.Range("A5").Value = "Drone With Camera Live Video, EACHINE E58 WIFI FPV Quadcopter With 120° Wide-angle 720P HD Camera Foldable Drone RTF - Altitude Hold, One Key Take Off/Landing, 3D Flip, APP Control, Gravity Sensor"
$ws.Range("B5").Value = "https://www.amazon.com/Camera-EACHINE-Quadcopter-Wide-angle-Foldable/dp/B0776QJNS3/ref=sr_1_7?s=toys-and-games&ie=UTF8&qid=1528131930&sr=1-7&keywords=drone"
$ws.Range("C5").Value = "4.2 out of 5 stars"
$ws.Range("D5").Value = 226
$ws.Range("D5").NumberFormat = "#,##0"
$ws.Range("E5").Value = 79

# Row 6
$ws.Range("A6").Value = "DROCON Drone For Beginners X708W Wi-Fi FPV Training Quadcopter With HD Camera Equipped With Headless Mode One Key Return Easy Operation"
$ws.Range("B6").Value = "https://www.amazon.com/DROCON-Beginners-Training-Quadcopter-Operation/dp/B073HYDPT3/ref=sr_1_8?s=toys-and-games&ie=UTF8&qid=1528131930&sr=1-8&keywords=drone"
$ws.Range("C6").Value = "4.2 out of 5 stars"
$ws.Range("D6").Value = 744
$ws.Range("E6").Value = 59

# Row 7
$ws.Range("C7").Value = "3.9 out of 5 stars"
$ws.Range("D7").Value = 4633
$ws.Range("E7").Value = 25

# Row 8
$ws.Range("A8").Value = "Cheerwing Syma X5SW-V3 FPV Explorers2 2.4Ghz 4CH 6-Axis Gyro RC Headless Quadcopter Drone UFO with HD Wifi Camera (White)"
$ws.Range("B8").Value = "https://www.amazon.com/Cheerwing-X5SW-V3-Explorers2-Headless-Quadcopter/dp/B011JV9HA2/ref=sr_1_10?s=toys-and-games&ie=UTF8&qid=1528131930&sr=1-10&keywords=drone"
$ws.Range("C8").Value = "3.7 out of 5 stars"
$ws.Range("D8").Value = 1513
$ws.Range("E8").Value = 39

# Row 9
$ws.Range("A9").Value = "Holy Stone F181C RC Quadcopter Drone with HD Camera RTF 4 Channel 2.4GHz 6-Gyro with Altitude Hold Function,Headless Mode and One Key Return Home, Color Black"
$ws.Range("B9").Value = "https://www.amazon.com/Holy-Stone-Quadcopter-Altitude-Function/dp/B00SAUAP5C/ref=sr_1_11?s=toys-and-games&ie=UTF8&qid=1528131930&sr=1-11&keywords=drone"
$ws.Range("C9").Value = "4.4 out of 5 stars"
$ws.Range("D9").Value = 3008
$ws.Range("E9").Value = 99

# Row 10
$ws.Range("A10").Value = "Mini Quadcopter Drone, EACHINE E010 2.4GHz 6-Axis Gyro Remote Control Nano Drone for Kids Adults Beginners - Headless Mode, 3D Flip, One Key Return (Green)"
$ws.Range("B10").Value = "https://www.amazon.com/Quadcopter-EACHINE-2-4GHz-Control-Beginners/dp/B01K9T54PC/ref=sr_1_12?s=toys-and-games&ie=UTF8&qid=1528131930&sr=1-12&keywords=drone"
$ws.Range("C10").Value = "3.8 out of 5 stars"
$ws.Range("D10").Value = 1698
$ws.Range("D10").NumberFormat = "#,##0"
$ws.Range("E10").Value = 24

# Row 11
$ws.Range("C11").Value = "4.2 out of 5 stars"
$ws.Range("D11").Value = 1175
$ws.Range("A11").ClearContents()
$ws.Range("B11").ClearContents()
$ws.Range("E11").ClearContents()
